$wb = $excel.ActiveWorkbook

# --- TextFileSequence sheet: remove the "name" column (G1) ---
$wsSeq = $wb.Worksheets.Item("TextFileSequence")
$wsSeq.Range("G1").ClearContents()

# --- Primer sheet: reorder header row so "name" comes first, ---
# --- shifting sequence/id/type one column to the right         ---
$wsPrimer = $wb.Worksheets.Item("Primer")
$wsPrimer.Range("A1").Value = "name"
$wsPrimer.Range("B1").Value = "sequence"
$wsPrimer.Range("C1").Value = "id"
$wsPrimer.Range("D1").Value = "type"
